# "update oops and metrics"
#
# The "Classrelationratio" metric row (row 14: Classrelationratio / Class
# Relation Ratio / ... / "Class/property ratio") is removed from the
# Foglio1 metrics table. Deleting the whole row shifts every subsequent
# row up by one (and the trailing blank-row styling pattern shifts with
# it), which reproduces the row 9/13-31 renumbering seen in the diff, and
# the now-unused shared strings are dropped automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("14").Delete()

# Restore the on-disk selection/view state: frozen top row, active cell on
# row 3 (the whole row selected), matching the target sheetView.
$ws.Rows("3:3").EntireRow.Select()
